$wb = $excel.ActiveWorkbook

# --- Sheet 1: "All Published Values" -----------------------------------
$ws1 = $wb.Worksheets.Item("All Published Values")

# Append the new data row (row 26) with the latest published BOC USD rate.
$rowNum = 26
$values = @(
    "2026-01-04",
    "2026-01-04 09:33:22",
    "697.35",
    "697.35",
    "700.29",
    "700.29",
    "702.88",
    "2026/01/04 09:33:22",
    "2026-01-04 01:58:58",
    "https://www.bankofchina.com/sourcedb/whpj/enindex_1619.html"
)

for ($col = 1; $col -le $values.Length; $col++) {
    $cell = $ws1.Cells.Item($rowNum, $col)
    # Force text storage so the date-looking strings aren't auto-converted
    # into date serial numbers, then reset the style back to the default
    # (unstyled) look used by the rest of the data rows.
    $cell.NumberFormat = "@"
    $cell.Value = $values[$col - 1]
    $cell.Style = "Normal"
}

# Re-apply the AutoFilter so its range grows to include the new row.
$ws1.AutoFilterMode = $false
$null = $ws1.Range("A1:J26").AutoFilter()

# Update the (hidden) _FilterDatabase defined name for this sheet so it
# also reflects the expanded range.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -eq "All Published Values!_FilterDatabase") {
        $n.RefersTo = "='All Published Values'!`$A`$1:`$J`$26"
    }
}

# --- Sheet 2: "Daily Summary" -------------------------------------------
$ws2 = $wb.Worksheets.Item("Daily Summary")

# The publish count for 2026-01-04 increases from 2 to 3 with the new entry.
$ws2.Range("B6").Value = 3
